# Apply "single meta-analysis and moderators" update to both
# sexmatch_strict_all_split (sheet1) and sexmatch_strict_all_split_0s (sheet2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update data values in column G (controls_NRH) / H (controls_RH) ---
# sheet1 (sexmatch_strict_all_split)
$ws1.Range("G2").Value = 128
$ws1.Range("H2").Value = 813
$ws1.Range("G3").Value = 127
$ws1.Range("H3").Value = 1089
$ws1.Range("G4").Value = 128
$ws1.Range("H4").Value = 813
$ws1.Range("G5").Value = 127
$ws1.Range("H5").Value = 1089
$ws1.Range("G25").Value = 4
$ws1.Range("H25").Value = 39
$ws1.Range("G27").Value = 4
$ws1.Range("H27").Value = 39

# sheet2 (sexmatch_strict_all_split_0s)
$ws2.Range("G2").Value = 128
$ws2.Range("H2").Value = 813
$ws2.Range("G3").Value = 127
$ws2.Range("H3").Value = 1089
$ws2.Range("G4").Value = 128
$ws2.Range("H4").Value = 813
$ws2.Range("G5").Value = 127
$ws2.Range("H5").Value = 1089
$ws2.Range("H25").Value = 39
$ws2.Range("H27").Value = 39

# --- Remove the total_cases / total_controls columns (I:J) from both sheets ---
$ws1.Columns("I:J").Delete()
$ws2.Columns("I:J").Delete()

# --- Restore selections: sheet1 gets F28 selected (not the active tab),
#     sheet2 keeps A20:H21 selected and remains the active tab. ---
$ws1.Range("F28").Select()
$ws2.Range("A20:H21").Select()
